$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty row 3 (B3:E3) with new test-case data,
# switching the retrieval method described in the commit message
# (API call replacing localStorage-only flow).
$ws.Range("B3").Value = "createSofa"
$ws.Range("C3").Value = "créer les nœuds HTML"
$ws.Range("D3").Value = "Afficher mes cards produits"
$ws.Range("E3").Value = "OK / Description erreur"

# Leave the selection parked the way the author saved the workbook.
$ws.Range("E4:E22").Select() | Out-Null
